$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix each hospital district location_code in column C (rows 3-23) with "fi_"
for ($r = 3; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = "fi_" + $cell.Value2
}

# Update the saved selection to match the edited range
$ws.Range("C3:C23").Select()
